$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "60.052.92"
$ws.Range("E2").Value = "  +1.35%  "

# Row 3
$ws.Range("D3").Value = "2.654.99"
$ws.Range("E3").Value = "  +2.47%  "

# Row 4
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "536.83"
$ws.Range("E5").Value = "  +1.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.95"
$ws.Range("E6").Value = "  +4.46%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.995"
$ws.Range("E7").Value = "  -0.31%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.574"
$ws.Range("E8").Value = "  +1.57%  "

# Row 9
$ws.Range("D9").Value = "2.672.20"
$ws.Range("E9").Value = "  +2.66%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.63"
$ws.Range("E10").Value = "  +2.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +2.65%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("E12").Value = "  +1.84%  "

# Row 13
$ws.Range("E13").Value = "  -1.69%  "

# Row 14
$ws.Range("D14").Value = "3.116.65"
$ws.Range("E14").Value = "  +2.19%  "

# Row 15
$ws.Range("D15").Value = "59.971.50"
$ws.Range("E15").Value = "  +1.40%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.17"
$ws.Range("E16").Value = "  +3.41%  "

# Row 17
$ws.Range("D17").Value = "2.661.48"
$ws.Range("E17").Value = "  +1.95%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000136"
$ws.Range("E18").Value = "  +1.80%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.62"
$ws.Range("E19").Value = "  -0.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.45"
$ws.Range("E20").Value = "  +2.79%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.30"
$ws.Range("E21").Value = "  +1.98%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.41"
$ws.Range("E22").Value = "  -0.60%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.80"
$ws.Range("E24").Value = "  +0.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.416"
$ws.Range("E25").Value = "  +2.56%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.167"
$ws.Range("E26").Value = "  +0.06%  "

# Row 27
$ws.Range("E27").Value = "  +0.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.33"
$ws.Range("E28").Value = "  +3.18%  "

# Row 29
$ws.Range("D29").Value = "0.0₃0759"
$ws.Range("E29").Value = "  +4.89%  "

# Row 30
$ws.Range("E30").Value = "  -0.04%  "

# Row 31
$ws.Range("E31").Value = "  +3.26%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.97"
$ws.Range("E32").Value = "  +1.27%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.10"
$ws.Range("E33").Value = "  +1.96%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.62"
$ws.Range("E34").Value = "  +1.13%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.05"
$ws.Range("E35").Value = "  +2.21%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.15"
$ws.Range("E36").Value = "  +2.96%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.48"
$ws.Range("E37").Value = "  +0.34%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.845"
$ws.Range("E38").Value = "  +2.36%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.834"
$ws.Range("E39").Value = "  +2.26%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "290.08"
$ws.Range("E40").Value = "  +7.66%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.62"
$ws.Range("E41").Value = "  +2.49%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.998"
$ws.Range("E42").Value = "  -0.12%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.607"
$ws.Range("E43").Value = "  +1.93%  "

# Row 44
$ws.Range("B44").Value = "Stellar"
$ws.Range("C44").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0961"
$ws.Range("E44").Value = "  +0.35%  "

# Row 45
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0541"
$ws.Range("E45").Value = "  +4.52%  "

# Row 46
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.74"
$ws.Range("E46").Value = "  -0.38%  "

# Row 47
$ws.Range("D47").Value = "1.978.11"
$ws.Range("E47").Value = "  +0.96%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0226"
$ws.Range("E48").Value = "  +2.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.60"
$ws.Range("E49").Value = "  +1.60%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.55"
$ws.Range("E50").Value = "  +1.77%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "110.61"
$ws.Range("E51").Value = "  -0.41%  "
